# Effort-Estimation date changes
# Each of the 4 weekly "Woche N until ..." headers gets a new date and a
# yellow highlight so the updated rows stand out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Woche 1 (row 4) -------------------------------------------------
$ws.Range("A4").Value = "Woche 1 until 21.03.2019"
$ws.Range("A4").Interior.Color = 65535
$ws.Range("A4").HorizontalAlignment = -4108

# --- Woche 2 (row 17) --------------------------------------------------
$ws.Range("A17").Value = "Woche 2 until 28.03.2019"
$ws.Range("A17").Interior.Color = 65535
$ws.Range("A17").HorizontalAlignment = -4108
$ws.Range("A17").VerticalAlignment = -4108

# --- Woche 3 (row 29) --------------------------------------------------
$ws.Range("A29").Value = "Woche 3 until 04.04.2019"
$ws.Range("A29").Interior.Color = 65535
$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("A29").VerticalAlignment = -4108

# --- Woche 4 (row 42) --------------------------------------------------
$ws.Range("A42").Value = "Woche 4 until 13.04.2019"
$ws.Range("A42").Interior.Color = 65535
$ws.Range("A42").HorizontalAlignment = -4108
$ws.Range("A42").VerticalAlignment = -4108
$ws.Range("A42").Font.Name = "Calibri"

# Move the active selection to C12 (matches the saved cursor position)
$ws.Range("C12").Select()
